# Refresh the per-coin "Price" (column D) and "Volume(1h)" (column E) text
# columns with this run's scraped values (scheduled GitHub Actions scrape).
#
# Both columns hold plain/General-formatted text (no numFmt on the cells) --
# e.g. D2 is the literal text "28.470.74" and E2 is "  +2.22%  " (note the
# padding spaces baked into the string). A handful of the new Price readings
# (like "315.10" or "41.94") are themselves valid numbers, so assigning them
# with a bare .Value would make Excel auto-convert the cell to a number and
# silently normalise away the trailing digit/zero (315.10 -> 315.1), which
# would not match the source text. Prefixing the literal with a leading
# apostrophe forces Excel to keep it as literal text -- exactly as if it had
# been typed into the cell by hand -- and then resetting Style back to
# "Normal" clears the quote-prefix flag Excel stamps on the cell so its style
# stays byte-identical to the original (plain default style, no numFmt).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextValue "D2" "28.470.74"
Set-TextValue "E2" "  +2.22%  "
Set-TextValue "D3" "1.828.64"
Set-TextValue "E3" "  +2.04%  "
Set-TextValue "E4" "  +0.13%  "
Set-TextValue "D5" "315.10"
Set-TextValue "E5" "  -0.38%  "
Set-TextValue "E6" "  +0.09%  "
Set-TextValue "D7" "0.5070"
Set-TextValue "E7" "  -4.70%  "
Set-TextValue "D8" "0.3908"
Set-TextValue "E8" "  +1.75%  "
Set-TextValue "D9" "0.07704"
Set-TextValue "E9" "  +3.81%  "
Set-TextValue "D10" "41.94"
Set-TextValue "E10" "  +1.48%  "
Set-TextValue "E11" "  +2.79%  "
Set-TextValue "D12" "21.02"
Set-TextValue "E12" "  +3.61%  "
Set-TextValue "D13" "6.268"
Set-TextValue "E13" "  +1.41%  "
Set-TextValue "D14" "7.571"
Set-TextValue "E14" "  +1.74%  "
Set-TextValue "E15" "  +0.12%  "
Set-TextValue "D16" "1.825.28"
Set-TextValue "E16" "  +1.90%  "
Set-TextValue "D17" "93.48"
Set-TextValue "E17" "  +6.20%  "
Set-TextValue "E18" "  +2.25%  "
Set-TextValue "D19" "0.06633"
Set-TextValue "E19" "  +1.78%  "
Set-TextValue "D20" "17.69"
Set-TextValue "E20" "  +2.69%  "
Set-TextValue "E21" "  +0.11%  "
Set-TextValue "E22" "  +3.50%  "
Set-TextValue "D23" "28.503.59"
Set-TextValue "E23" "  +2.19%  "
Set-TextValue "D24" "11.15"
Set-TextValue "E24" "  +0.27%  "
Set-TextValue "D25" "2.256"
Set-TextValue "E25" "  +7.92%  "
Set-TextValue "D26" "156.78"
Set-TextValue "E26" "  -0.20%  "
Set-TextValue "D27" "20.63"
Set-TextValue "E27" "  +2.63%  "
Set-TextValue "D28" "2.036.94"
Set-TextValue "E28" "  +1.98%  "
Set-TextValue "D29" "2.419"
Set-TextValue "E29" "  +4.38%  "
Set-TextValue "D30" "125.61"
Set-TextValue "E30" "  +3.68%  "
Set-TextValue "E31" "  +3.34%  "
Set-TextValue "D32" "0.1089"
Set-TextValue "E32" "  -0.21%  "
Set-TextValue "D33" "5.660"
Set-TextValue "E33" "  +2.98%  "
Set-TextValue "D34" "3.658"
Set-TextValue "E34" "  +0.29%  "
Set-TextValue "D35" "0.07058"
Set-TextValue "E35" "  +2.12%  "
Set-TextValue "D36" "0.2219"
Set-TextValue "E36" "  +0.82%  "
Set-TextValue "D37" "8.981"
Set-TextValue "E37" "  +7.42%  "
Set-TextValue "D38" "0.02328"
Set-TextValue "E38" "  +2.63%  "
Set-TextValue "D39" "5.156"
Set-TextValue "E39" "  +2.47%  "
Set-TextValue "D40" "0.6246"
Set-TextValue "E40" "  +2.51%  "
Set-TextValue "D41" "11.19"
Set-TextValue "E41" "  -1.03%  "
Set-TextValue "E42" "  +1.40%  "
Set-TextValue "E43" "  +0.03%  "
Set-TextValue "D44" "1.397"
Set-TextValue "E44" "  -0.90%  "
Set-TextValue "D45" "13.45"
Set-TextValue "E45" "  +1.39%  "
Set-TextValue "D46" "0.5905"
Set-TextValue "E46" "  +3.78%  "
Set-TextValue "D47" "3.715"
Set-TextValue "E47" "  +0.96%  "
Set-TextValue "D48" "124.77"
Set-TextValue "E48" "  +0.35%  "
Set-TextValue "D49" "1.975"
Set-TextValue "E49" "  +3.43%  "
Set-TextValue "D50" "1.200"
Set-TextValue "E50" "  +2.67%  "
Set-TextValue "D51" "0.06927"
Set-TextValue "E51" "  +2.00%  "
